$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert a new row at position 20 (shifts existing rows 20-22 down to 21-23)
$ws.Range("A20").EntireRow.Insert()

# 2) Populate the new row 20 with the outagesFetchUrl entry
$ws.Range("A20").Value = "outagesFetchUrl"
$ws.Range("B20").Value = "http://google.com"
# copy format (style) from a neighboring hyperlink-styled cell so B20 gets style index 1 (Hyperlink)
$ws.Range("B19").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) The pre-existing hyperlink relationship that used to point at B20 (transOutagesFetchUrl row)
#    now naturally (and conveniently) decorates the new B20 cell, since its stored range reference
#    was never shifted by the row insert. Re-assert it explicitly (exact-range match => updates
#    the existing relationship in place rather than creating a duplicate).
$ws.Hyperlinks.Add($ws.Range("B20"), "http://google.com/") | Out-Null
$ws.Range("B19").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4) Remove the now-incorrect merged hyperlink that used to span B21:B22 (majorGen + longUnrev).
#    This relationship refers to a multi-cell range and can be cleanly deleted.
$ws.Hyperlinks.Item(14).Delete()

# 5) Re-create the hyperlink for the merged majorGen/longUnrev rows at their new location B22:B23
$ws.Hyperlinks.Add($ws.Range("B22:B23"), "http://google.com/", [Type]::Missing, "http://google.com") | Out-Null

# 6) Give the transOutagesFetchUrl row (now row 21) its own hyperlink, since the original one
#    stayed behind on B20
$ws.Hyperlinks.Add($ws.Range("B21"), "http://google.com/") | Out-Null

# 7) Match the recorded selection from the authored workbook
$ws.Range("A21").Select()

Write-Host "done"
